# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1016
$ws1.Range("F4").Value  = 13402
$ws1.Range("F9").Value  = 129
$ws1.Range("F10").Value = 117
$ws1.Range("F14").Value = 13379
$ws1.Range("F17").Value = 8917
$ws1.Range("F18").Value = 3
$ws1.Range("F19").Value = 7984
$ws1.Range("F20").Value = 244
$ws1.Range("F21").Value = 6
$ws1.Range("F32").Value = 161
$ws1.Range("F33").Value = 370

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1016
$ws4.Range("F5").Value  = 13402
$ws4.Range("F10").Value = 129
$ws4.Range("F11").Value = 117
$ws4.Range("F15").Value = 13379
$ws4.Range("F18").Value = 8917
$ws4.Range("F19").Value = 3
$ws4.Range("F20").Value = 7984
$ws4.Range("F21").Value = 244
$ws4.Range("F22").Value = 6
$ws4.Range("F35").Value = 161
$ws4.Range("F36").Value = 370
